# Fix the LOQ4243 syllabus sheet:
#  - insert a row for "Docentes responsáveis:" value (previously missing)
#  - correct several mis-shifted / placeholder cell values that had been
#    copy-pasted from the wrong rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 13 (shifts "Programa resumido:" and everything below
# down by one row, carrying their row heights with them).
$ws.Rows(13).Insert()

# Objetivos (row 10) had the wrong value copied in - put the real
# Portuguese objectives text in B10/C10.
$ws.Range("B10").Value = "Este curso visa apresentar os principais conceitos relacionados ao uso da eletrônica, com foco na identificação de componentes e funcionalidades, para o uso em projetos relacionados ao curso de graduação em engenharia de produção."
$ws.Range("C10").Value = "Este curso visa apresentar os principais conceitos relacionados ao uso da eletrônica, com foco na identificação de componentes e funcionalidades, para o uso em projetos relacionados ao curso de graduação em engenharia de produção."

# Newly inserted row 13 holds the "Docentes responsáveis:" value.
$ws.Range("B13").Value = "8767640 - Eduardo Ferro dos Santos"
$ws.Range("C13").Value = "8767640 - Eduardo Ferro dos Santos"
# The blank inserted row leaves A13 with an inherited border/format and
# B13 with the wrong (default) style - clean those up to match the rest
# of the table (no A-cell on this row, B uses the same wrap style as C).
$ws.Range("A13").Clear()
$ws.Range("B14").Copy()
$ws.Range("B13").PasteSpecial(-4122)

# Programa resumido (now row 14) had "Semestral" left over - replace with
# the actual short syllabus text.
$ws.Range("B14").Value = "Fundamentos da Eletrônica; Sensores e atuadores; Microcontroladores;Sistemas mecatrônicos;Experimentos práticos."
$ws.Range("C14").Value = "Fundamentos da Eletrônica; Sensores e atuadores; Microcontroladores;Sistemas mecatrônicos;Experimentos práticos."

# Programa (now row 16) had a stray date value - replace with the full
# syllabus text.
$ws.Range("B16").Value = "Fundamentos da Eletrônica: digital e analógica; Uso de sensores e atuadores em processos produtivos; Programação básica de microcontroladores, com foco em Arduino ou similar;Projetos de uso em sistemas mecatrônicos aplicados a engenharia de produção;Experimentos práticos."
$ws.Range("C16").Value = "Fundamentos da Eletrônica: digital e analógica; Uso de sensores e atuadores em processos produtivos; Programação básica de microcontroladores, com foco em Arduino ou similar;Projetos de uso em sistemas mecatrônicos aplicados a engenharia de produção;Experimentos práticos."

# Método (now row 19) had the docente name leftover - replace with the
# evaluation method text.
$ws.Range("B19").Value = "Esta disciplina deverá conter no mínimo duas avaliações denominadas A1 e A2. As avalições poderão ser: escritas, práticas, seminários, trabalhos de campo, projetos, ou outra forma de avaliação definida pelo professor."
$ws.Range("C19").Value = "Esta disciplina deverá conter no mínimo duas avaliações denominadas A1 e A2. As avalições poderão ser: escritas, práticas, seminários, trabalhos de campo, projetos, ou outra forma de avaliação definida pelo professor."

# Critério (now row 20) should show the weighted-average criterion text.
$ws.Range("B20").Value = "Média ponderada das avaliações (M)."
$ws.Range("C20").Value = "Média ponderada das avaliações (M)."

# Norma de recuperação (now row 21) should show the recovery-exam text.
$ws.Range("B21").Value = "A recuperação será composta por uma única prova (RC) englobando toda a matéria ministrada ao longo do semestre. A média final, para os alunos em recuperação, será calculada com base na relação: MF=(M+RC)/2"
$ws.Range("C21").Value = "A recuperação será composta por uma única prova (RC) englobando toda a matéria ministrada ao longo do semestre. A média final, para os alunos em recuperação, será calculada com base na relação: MF=(M+RC)/2"

# Bibliografia (now row 22) should show the actual bibliography text.
$ws.Range("B22").Value = "Tutoriais de Arduino disponibilizados pelo fabricante (arduino.cc) BOYLESTAD, Robert L.; NASHELSKY, Louis. Dispositivos Eletrônicos e Teoria de Circuitos. 8ª ed. São Paulo: Pearson. 696 p. THOMAZINI, Daniel; ALBUQUERQUE, Pedro U.B. Sensores Industriais – Fundamentos e Aplicações. 8ª ed. São Paulo: Érica, 2011. 224 p.Bibliografia complementar será indicada ao longo do curso."
$ws.Range("C22").Value = "Tutoriais de Arduino disponibilizados pelo fabricante (arduino.cc) BOYLESTAD, Robert L.; NASHELSKY, Louis. Dispositivos Eletrônicos e Teoria de Circuitos. 8ª ed. São Paulo: Pearson. 696 p. THOMAZINI, Daniel; ALBUQUERQUE, Pedro U.B. Sensores Industriais – Fundamentos e Aplicações. 8ª ed. São Paulo: Érica, 2011. 224 p.Bibliografia complementar será indicada ao longo do curso."
